$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2603.423076923077, 44260.53459715942),
    @(2544.222222222222, 44260.54520345246),
    @(2547.259259259259, 44260.54536520431),
    @(7, 44260.54552288155),
    @(2550.148148148148, 44260.54738223815),
    @(2, 44260.54751858036),
    @(68768.55555555556, 44260.55013614835),
    @(2549.259259259259, 44260.57246532349),
    @(-2521.814814814815, 44260.57480656946),
    @(3446.444444444444, 44260.57951437631),
    @(2549.074074074074, 44260.57973993074),
    @(6914.333333333333, 44260.63564280069),
    @(-1.148148148148148, 44260.63583176896),
    @(2549.62962962963, 44260.65864267803),
    @(19495.14814814815, 44260.73166254783),
    @(2548.481481481481, 44260.73209902812),
    @(7505.089251804026, 44260.73854107285),
    @(15005.84884162552, 44260.7390276044),
    @(2617.242688947968, 44260.74161344661)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("B2:B20").NumberFormat = "yyyy-mm-dd h:mm:ss"
